# Auto-generated Excel COM-interop script to apply cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.602.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -10.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.358.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -13.68%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "453.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -10.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.475"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -11.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.369.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -13.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0929"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -11.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.17"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -16.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.308"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -12.02%  "
$ws.Range("E13").Value = "  -4.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.766.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -13.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "52.706.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -10.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -11.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000127"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.370.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -12.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -14.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "302.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -13.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -16.62%  "
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -16.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "54.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -12.76%  "
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.376"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -11.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.476.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -12.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.148"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -13.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.83%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0708"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -15.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "143.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.56%  "
$ws.Range("E35").Value = "  -13.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -18.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.68%  "
$ws.Range("E39").Value = "  -18.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.993"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "32.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.581"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.65%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0514"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.32%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("E46").Value = "  -13.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.915.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -12.69%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0213"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.45%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0849"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.39%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -14.91%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -15.76%  "
